$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 9 (August) label to reflect new "through" date
$ws.Range("A9").Value = "August (through 08-28)"

# Update row 9 values (B9 unchanged)
$ws.Range("C9").Value = 68
$ws.Range("D9").Value = 82
$ws.Range("E9").Value = 59
$ws.Range("F9").Value = 42
$ws.Range("G9").Value = 154
$ws.Range("H9").Value = 144

# Update row 10 (Total) values (B10 unchanged)
$ws.Range("C10").Value = 370
$ws.Range("D10").Value = 547
$ws.Range("E10").Value = 484
$ws.Range("F10").Value = 346
$ws.Range("G10").Value = 775
$ws.Range("H10").Value = 1058
